$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 328
$ws.Range("F4").Value = 3011
$ws.Range("F7").Value = 2340
$ws.Range("F8").Value = 1730
$ws.Range("F9").Value = 63
$ws.Range("F11").Value = 137
$ws.Range("F13").Value = 25
$ws.Range("F14").Value = 2678
$ws.Range("F16").Value = 1551
$ws.Range("F17").Value = 7172
$ws.Range("F19").Value = 7311
$ws.Range("F22").Value = 5619
$ws.Range("F23").Value = 3138
$ws.Range("F24").Value = 3515
$ws.Range("F25").Value = 2
$ws.Range("F26").Value = 7
$ws.Range("F27").Value = 247
$ws.Range("F28").Value = 198
$ws.Range("F29").Value = 1943
$ws.Range("F31").Value = 312
$ws.Range("F33").Value = 230
$ws.Range("F34").Value = 497
$ws.Range("F36").Value = 2466
$ws.Range("F37").Value = 1256
$ws.Range("F38").Value = 2843
$ws.Range("F39").Value = 60
$ws.Range("F41").Value = 176
$ws.Range("F43").Value = 1116
$ws.Range("F45").Value = 491
$ws.Range("F46").Value = 541

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 41
$ws.Range("F10").Value = 10
$ws.Range("F12").Value = 354
$ws.Range("F13").Value = 25
$ws.Range("F14").Value = 50
$ws.Range("F15").Value = 97
$ws.Range("F19").Value = 70
$ws.Range("F21").Value = 8

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 87

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 328
$ws.Range("F5").Value = 3011
$ws.Range("F6").Value = 2340
$ws.Range("F7").Value = 1730
$ws.Range("F8").Value = 63
$ws.Range("F10").Value = 137
$ws.Range("F12").Value = 41
$ws.Range("F13").Value = 87
$ws.Range("F14").Value = 2678
$ws.Range("F15").Value = 1551
$ws.Range("F19").Value = 7172
$ws.Range("F21").Value = 7311
$ws.Range("F23").Value = 5619
$ws.Range("F24").Value = 3138
$ws.Range("F25").Value = 3515
$ws.Range("F26").Value = 7
$ws.Range("F27").Value = 25
$ws.Range("F28").Value = 247
$ws.Range("F29").Value = 50
$ws.Range("F30").Value = 1943
$ws.Range("F33").Value = 312
$ws.Range("F35").Value = 230
$ws.Range("F36").Value = 497
$ws.Range("F38").Value = 2466
$ws.Range("F39").Value = 1256
$ws.Range("F40").Value = 70
$ws.Range("F41").Value = 2843
$ws.Range("F42").Value = 60
$ws.Range("F44").Value = 176
$ws.Range("F46").Value = 1116
$ws.Range("F48").Value = 491
$ws.Range("F49").Value = 541
